$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.129.35"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "1.680.43"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'215.20"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "'21.33"
$ws.Range("E9").Value = "  +5.29%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.0623"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "1.917.20"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "1.688.25"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "'66.21"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "27.115.03"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "'238.19"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").Value = "'8.13"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'4.52"
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("D23").Value = "'9.48"
$ws.Range("E23").Value = "  +3.32%  "
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("D25").Value = "'146.80"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").Value = "'7.24"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").Value = "'16.31"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").Value = "1.564.07"
$ws.Range("E32").Value = "  +5.53%  "
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("E36").Value = "  +3.18%  "
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("D38").Value = "'0.937"
$ws.Range("E38").Value = "  +4.44%  "
$ws.Range("D39").Value = "'0.0175"
$ws.Range("E39").Value = "  +2.63%  "
$ws.Range("D40").Value = "'1.06"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("D41").Value = "'69.04"
$ws.Range("E41").Value = "  +2.81%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "'5.66"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("E44").Value = "  -1.81%  "
$ws.Range("D45").Value = "1.825.34"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'90.78"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  +3.20%  "
$ws.Range("E49").Value = "  +1.69%  "
$ws.Range("D50").Value = "'0.105"
$ws.Range("E50").Value = "  +2.69%  "
$ws.Range("D51").Value = "'8.09"
$ws.Range("E51").Value = "  +5.02%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
